$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-06 01:17:09"

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
